# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted at row 17 (pushing the
# previously-existing rows 17-90 down to 18-91); the new row holds the
# most recent observation (2021-10-18) and the rest of the table is
# otherwise unchanged, just shifted down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 17. Excel shifts every
# row from 17 downward (previously 17..90) down to 18..91, and widens
# the used range / dimension to A1:R91 automatically.
$ws.Rows(17).Insert()

# Populate the newly inserted row 17 with the new weekly observation.
$ws.Cells.Item(17, 1).Value2 = 7
$ws.Cells.Item(17, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value2 = "Ñuble"
$ws.Cells.Item(17, 4).Value2 = 44487
$ws.Cells.Item(17, 5).Value2 = 16
$ws.Cells.Item(17, 6).Value2 = 100112045
$ws.Cells.Item(17, 7).Value2 = "Zapallo"
$ws.Cells.Item(17, 8).Value2 = "Camote"
$ws.Cells.Item(17, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(17, 10).Value2 = 120
$ws.Cells.Item(17, 11).Value2 = 800
$ws.Cells.Item(17, 12).Value2 = 900
$ws.Cells.Item(17, 13).Value2 = 850
$ws.Cells.Item(17, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(17, 15).Value2 = "Región del Maule"
$ws.Cells.Item(17, 16).Value2 = 850
$ws.Cells.Item(17, 17).Value2 = 1
$ws.Cells.Item(17, 18).Value2 = "Hortaliza"
